# Fiscal year code field added
#
# The "fiscalyearmaster" table (row 6) and the "functiondetails" table
# (row 7) on the "Tables" sheet both gain a leading "code" field.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Fields" column (B) for the functiondetails row (row 7) first,
# then the fiscalyearmaster row (row 6), so new shared-string entries are
# appended in the same order as the reference edit.
$ws.Range("B7").Value = "code,format,digit_length"
$ws.Range("B6").Value = "code,name,abr,start_date,end_date"

# Move the active selection to B6, matching the reference edit.
$ws.Range("B6").Select() | Out-Null
